$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "I" = 0.01686872651965528
    "J" = 0.01402212891946345
    "K" = 55.4
    "L" = 0.3325330132052821
    "U" = 443.4
    "V" = 0.148076409297355
    "W" = 0.4140508221225709
    "X" = 0.03990200775101115
    "Y" = 0.3741488143715598
    "Z" = 0.9970772899082292
    "AA" = 0.01398114630176242
    "AB" = 0.03917168479431483
    "AC" = -0.02519053849255241
    "AD" = 249.3
    "AE" = 15.79835080912715
    "AF" = 265.0983508091272
    "AG" = -178.3016491908728
    "AH" = 0.08133102774643836
    "AI" = 0.6014050422886429
    "AJ" = -0.06331513568751702
    "AK" = 68.53408592380393
    "AN" = 41.75879396984925
    "AP" = -29.86627289629361
}

foreach ($col in $updates.Keys) {
    $value = $updates[$col]
    $ws.Range("${col}2").Value = $value
    $ws.Range("${col}3").Value = $value
}
